# Recategorize "Designation" workbook:
#  - rename a few category labels in column D (Groupe) that were judged
#    too verbose / redundant
#  - recategorize one row (Code de l'action sociale et des familles) from
#    the (now retired) "Loi diverses" bucket into "Admin divers"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename categories (these apply to every cell using the old label,
# and Excel will drop the now-unused shared string automatically)
$ws.Cells.Replace("Science et littérature", "Science, littérature")
$ws.Cells.Replace("Organisations internationales", "Org Internationales")
$ws.Cells.Replace("Expertise, Forums, Autorités", "Forums, Autorités")

# Recategorize the single remaining "Loi diverses" row
$ws.Range("C173").Value = "Admin divers"

# Restore the view/selection state recorded the last time the workbook was saved
$ws.Activate()
$ws.Range("I50").Select()

Write-Host "recategorization applied"
